$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.962.94'
$ws.Range("E2").Value = '  +0.65%  '
$ws.Range("D3").Value = '1.554.83'
$ws.Range("E3").Value = '  +0.33%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.01'
$ws.Range("E4").Value = '  +0.44%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '206.81'
$ws.Range("E5").Value = '  +0.76%  '
$ws.Range("E6").Value = '  +0.37%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.41%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '21.67'
$ws.Range("E8").Value = '  +1.16%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.248'
$ws.Range("E9").Value = '  +0.87%  '
$ws.Range("E10").Value = '  +0.47%  '
$ws.Range("E11").Value = '  +0.16%  '
$ws.Range("D12").Value = '1.775.15'
$ws.Range("E12").Value = '  +0.32%  '
$ws.Range("D13").Value = '1.552.88'
$ws.Range("E13").Value = '  +0.22%  '
$ws.Range("E14").Value = '  +0.71%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.516'
$ws.Range("E15").Value = '  +0.86%  '
$ws.Range("D16").Value = '26.957.23'
$ws.Range("E16").Value = '  +0.71%  '
$ws.Range("E17").Value = '  +1.42%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '214.48'
$ws.Range("E18").Value = '  +0.32%  '
$ws.Range("D19").Value = '0.0₃0687'
$ws.Range("E19").Value = '  +0.13%  '
$ws.Range("E20").Value = '  +0.33%  '
$ws.Range("E21").Value = '  +0.39%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.04'
$ws.Range("E22").Value = '  -1.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.19'
$ws.Range("E23").Value = '  +1.85%  '
$ws.Range("E24").Value = '  -1.30%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.14'
$ws.Range("E25").Value = '  +0.22%  '
$ws.Range("E26").Value = '  +2.28%  '
$ws.Range("E27").Value = '  -0.05%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.01'
$ws.Range("E28").Value = '  +0.44%  '
$ws.Range("E29").Value = '  +1.25%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0461'
$ws.Range("E30").Value = '  -0.08%  '
$ws.Range("E31").Value = '  -0.69%  '
$ws.Range("E32").Value = '  +1.69%  '
$ws.Range("D33").Value = '1.379.08'
$ws.Range("E33").Value = '  +2.01%  '
$ws.Range("E34").Value = '  +2.54%  '
$ws.Range("E35").Value = '  +3.47%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.974'
$ws.Range("E36").Value = '  +6.09%  '
$ws.Range("E37").Value = '  +0.68%  '
$ws.Range("E38").Value = '  +1.66%  '
$ws.Range("E39").Value = '  -0.13%  '
$ws.Range("E40").Value = '  +0.84%  '
$ws.Range("E41").Value = '  +0.42%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.990'
$ws.Range("E42").Value = '  -0.08%  '
$ws.Range("B43").Value = 'MXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.26'
$ws.Range("E43").Value = '  +3.15%  '
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.47'
$ws.Range("E44").Value = '  -1.71%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '63.90'
$ws.Range("E45").Value = '  +1.54%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.74'
$ws.Range("E46").Value = '  -1.40%  '
$ws.Range("D47").Value = '1.688.40'
$ws.Range("E47").Value = '  +0.34%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '86.28'
$ws.Range("E48").Value = '  +0.45%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0510'
$ws.Range("E49").Value = '  +0.66%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0954'
$ws.Range("E50").Value = '  +1.09%  '
$ws.Range("E51").Value = '  +0.55%  '
